# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before "总计", populate it with
#    the fund-holdings detail for that quarter (same layout/style as the
#    other quarterly sheets).
# 2. Insert a new row at the top of the data in "总计" for "2022-Q1" and
#    push the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("2021-Q4")
$totalBeforeAdd = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. New "2022-Q1" detail sheet, inserted immediately before "总计".
#    NOTE: sheet object handles in this host resolve by tab position, so
#    once the new sheet is inserted, the handle captured above
#    ($totalBeforeAdd) now refers to whatever sits at that old index
#    (i.e. the new sheet itself) - it must not be used afterwards. Grab
#    the freshly created sheet's own handle instead, and re-look-up
#    "总计" by name for part 2 below.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($totalBeforeAdd)
$q1.Name = "2022-Q1"

# Pull header-row / index-column formatting (bold, centered, bordered)
# from the existing "2021-Q4" sheet so the new sheet matches the others.
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$template.Range("A2:A8").Copy()
$q1.Range("A2:A8").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$fundRows = @(
    @{ Row=2; A=0; Code="519778"; Name="交银施罗德经济新动力混合"; Scale="45.93"; Pos="63.82"; Pct="5.55"; Mv="2.5491"; Rank=4 },
    @{ Row=3; A=1; Code="008955"; Name="交银施罗德创新领航混合";   Scale="32.34"; Pos="65.83"; Pct="5.60"; Mv="1.8110"; Rank=4 },
    @{ Row=4; A=2; Code="013247"; Name="交银瑞卓三年持有期混合";   Scale="19.29"; Pos="60.51"; Pct="5.48"; Mv="1.0571"; Rank=4 },
    @{ Row=5; A=3; Code="519727"; Name="交银成长30混合";           Scale="12.41"; Pos="64.16"; Pct="5.62"; Mv="0.6974"; Rank=4 },
    @{ Row=6; A=4; Code="003292"; Name="嘉实优势成长灵活配置混合"; Scale="1.04";  Pos="92.17"; Pct="4.93"; Mv="0.0513"; Rank=9 },
    @{ Row=7; A=5; Code="002630"; Name="江信瑞福灵活配置混合A";    Scale="0.52";  Pos="43.17"; Pct="2.50"; Mv="0.0130"; Rank=7 },
    @{ Row=8; A=6; Code="002631"; Name="江信瑞福灵活配置混合C";    Scale="0.50";  Pos="43.17"; Pct="2.50"; Mv="0.0125"; Rank=7 }
)

foreach ($fr in $fundRows) {
    $r = $fr.Row

    $q1.Range("A$r").Value = $fr.A

    # Columns that look numeric (fund code, scale, position, pct, market
    # value) must stay text, matching the source data, so force a text
    # number-format before writing, then drop back to the default style
    # afterwards (keeps the "text" cell type but leaves no numberformat
    # override behind, matching the un-styled cells in the reference diff).
    $codeCell = $q1.Range("B$r")
    $numericTextRange = $q1.Range("D$r:G$r")
    $codeCell.NumberFormat = "@"
    $numericTextRange.NumberFormat = "@"

    $q1.Range("B$r").Value = $fr.Code
    $q1.Range("C$r").Value = $fr.Name
    $q1.Range("D$r").Value = $fr.Scale
    $q1.Range("E$r").Value = $fr.Pos
    $q1.Range("F$r").Value = $fr.Pct
    $q1.Range("G$r").Value = $fr.Mv

    $codeCell.Style = "Normal"
    $numericTextRange.Style = "Normal"

    $q1.Range("H$r").Value = $fr.Rank
}

# ---------------------------------------------------------------------
# 2. "总计" sheet: insert a new "2022-Q1" row above the existing data.
#    Re-look-up the sheet by name now that the tab order has changed.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$dateCell = $total.Range("B2")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2022-Q1"
$dateCell.Style = "Normal"

$total.Range("C2").Style = "Normal"
$total.Range("D2").Style = "Normal"

$total.Range("A2").Value = 0
$total.Range("C2").Value = 7
$total.Range("D2").Value = 6.19

# Restore the bold/bordered index-column style on the new row's A cell
# (the row-insert above left it with the default/unstyled format).
$template.Range("A2").Copy()
$total.Range("A2").PasteSpecial(-4122)
